# Apply cryptos-list price/volume refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.106.60"
$ws.Range("E2").Value = "'  -0.77%  "
$ws.Range("D3").Value = "'1.652.30"
$ws.Range("E3").Value = "'  -0.88%  "
$ws.Range("E4").Value = "'  -0.44%  "
$ws.Range("D5").Value = "'218.77"
$ws.Range("E5").Value = "'  -0.77%  "
$ws.Range("D6").Value = "'0.5254"
$ws.Range("E6").Value = "'  -1.01%  "
$ws.Range("E7").Value = "'  -0.45%  "
$ws.Range("D8").Value = "'0.2673"
$ws.Range("E8").Value = "'  +1.10%  "
$ws.Range("D9").Value = "'0.06369"
$ws.Range("E9").Value = "'  +0.17%  "
$ws.Range("D10").Value = "'20.55"
$ws.Range("E10").Value = "'  -1.70%  "
$ws.Range("D11").Value = "'0.07681"
$ws.Range("E11").Value = "'  -1.97%  "
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.731.11"
$ws.Range("E12").Value = "'  +3.74%  "
$ws.Range("B13").Value = "'Polkadot"
$ws.Range("C13").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.596"
$ws.Range("E13").Value = "'  +1.48%  "
$ws.Range("D14").Value = "'1.880.27"
$ws.Range("E14").Value = "'  -0.80%  "
$ws.Range("D15").Value = "'0.5606"
$ws.Range("E15").Value = "'  +0.02%  "
$ws.Range("D16").Value = "'0.0₅8236"
$ws.Range("E16").Value = "'  +1.31%  "
$ws.Range("E17").Value = "'  -0.48%  "
$ws.Range("D18").Value = "'26.114.82"
$ws.Range("E18").Value = "'  -0.77%  "
$ws.Range("E19").Value = "'  -0.40%  "
$ws.Range("D20").Value = "'4.691"
$ws.Range("E20").Value = "'  -0.39%  "
$ws.Range("D21").Value = "'10.36"
$ws.Range("E21").Value = "'  +0.87%  "
$ws.Range("D22").Value = "'191.25"
$ws.Range("E22").Value = "'  -3.92%  "
$ws.Range("D23").Value = "'5.975"
$ws.Range("E23").Value = "'  -1.25%  "
$ws.Range("E24").Value = "'  -0.51%  "
$ws.Range("D25").Value = "'146.01"
$ws.Range("E25").Value = "'  -0.42%  "
$ws.Range("D26").Value = "'0.1201"
$ws.Range("E26").Value = "'  -0.95%  "
$ws.Range("D27").Value = "'7.254"
$ws.Range("E27").Value = "'  +0.33%  "
$ws.Range("E28").Value = "'  -1.22%  "
$ws.Range("E29").Value = "'  -1.00%  "
$ws.Range("D30").Value = "'0.05646"
$ws.Range("E30").Value = "'  -4.16%  "
$ws.Range("E31").Value = "'  -0.95%  "
$ws.Range("D32").Value = "'3.501"
$ws.Range("E32").Value = "'  -1.01%  "
$ws.Range("D33").Value = "'3.380"
$ws.Range("E33").Value = "'  +1.98%  "
$ws.Range("D34").Value = "'1.580"
$ws.Range("E34").Value = "'  -1.30%  "
$ws.Range("E35").Value = "'  -1.03%  "
$ws.Range("D36").Value = "'0.9453"
$ws.Range("E36").Value = "'  -1.53%  "
$ws.Range("D37").Value = "'2.407"
$ws.Range("E37").Value = "'  -0.96%  "
$ws.Range("D38").Value = "'0.5782"
$ws.Range("E38").Value = "'  -0.33%  "
$ws.Range("D39").Value = "'0.01592"
$ws.Range("E39").Value = "'  -1.58%  "
$ws.Range("D40").Value = "'5.974"
$ws.Range("E40").Value = "'  +0.20%  "
$ws.Range("E41").Value = "'  -0.53%  "
$ws.Range("D42").Value = "'0.8408"
$ws.Range("E42").Value = "'  -1.84%  "
$ws.Range("D43").Value = "'1.025.07"
$ws.Range("E43").Value = "'  -4.54%  "
$ws.Range("D44").Value = "'101.67"
$ws.Range("E44").Value = "'  -1.06%  "
$ws.Range("D45").Value = "'1.791.93"
$ws.Range("E45").Value = "'  -0.78%  "
$ws.Range("D46").Value = "'58.57"
$ws.Range("E46").Value = "'  +0.32%  "
$ws.Range("E47").Value = "'  -0.89%  "
$ws.Range("D48").Value = "'0.05339"
$ws.Range("E48").Value = "'  +3.74%  "
$ws.Range("B49").Value = "'BabyDogeCoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.0₈103"
$ws.Range("E49").Value = "'  -1.12%  "
$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.025"
$ws.Range("E50").Value = "'  -0.16%  "
$ws.Range("D51").Value = "'0.4343"
